# Updated remaining queries for C3DC
# - Rewrite the join conditions in every SQL query cell (C2, B2:B7) so that
#   they use the renamed `study_id` / `participant_id` columns instead of
#   the old generic `id` columns.
# - Move the active selection to B2 (instead of C5) and scroll the sheet
#   back to the top.
# - Widen column C now that it no longer relies on "best fit".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-Query($cellRef) {
    $cell = $ws.Range($cellRef)
    $text = $cell.Value2
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $cell.Value = $text
}

# StatQuery (C2) + TabQuery cells (B2:B7) all share the same stale join
# conditions and need the same fix.
Update-Query "C2"
Update-Query "B2"
Update-Query "B3"
Update-Query "B4"
Update-Query "B5"
Update-Query "B6"
Update-Query "B7"

# Column C no longer needs to rely on auto "best fit" sizing - widen it
# explicitly.
$ws.Range("C1").ColumnWidth = 67.3

# Selection moves to B2 and the view scrolls back to show row 1.
[void]$ws.Range("B2").Select()
